$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = "NK Domzale"
# Row 5
$ws.Cells.Item(5, 7).Value = "NK Maribor"
# Row 9
$ws.Cells.Item(9, 2).Value = 6814328
$ws.Cells.Item(9, 7).Value = "NK Bravo"
$ws.Cells.Item(9, 9).Value = 1
$ws.Cells.Item(9, 10).Value = "D"
$ws.Cells.Item(9, 11).Value = 2.35
$ws.Cells.Item(9, 12).Value = 3.1
$ws.Cells.Item(9, 13).Value = 2.9
$ws.Cells.Item(9, 14).Value = 2.15
$ws.Cells.Item(9, 15).Value = 3.1
$ws.Cells.Item(9, 16).Value = 3.3
$ws.Cells.Item(9, 17).Value = -0.25
$ws.Cells.Item(9, 18).Value = 1.925
$ws.Cells.Item(9, 19).Value = 1.875
$ws.Cells.Item(9, 20).Value = 2.25
$ws.Cells.Item(9, 21).Value = 1.95
$ws.Cells.Item(9, 22).Value = 1.85
$ws.Cells.Item(9, 23).Value = -1
$ws.Cells.Item(9, 24).Value = 2.1
$ws.Cells.Item(9, 27).Value = 0.4375
$ws.Cells.Item(9, 28).Value = -0.5
$ws.Cells.Item(9, 29).Value = 0.425
# Row 10
$ws.Cells.Item(10, 2).Value = 6814330
$ws.Cells.Item(10, 7).Value = "NK Aluminij"
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = "H"
$ws.Cells.Item(10, 11).Value = 1.363
$ws.Cells.Item(10, 12).Value = 4.5
$ws.Cells.Item(10, 13).Value = 7
$ws.Cells.Item(10, 14).Value = 1.4
$ws.Cells.Item(10, 15).Value = 4.5
$ws.Cells.Item(10, 16).Value = 7
$ws.Cells.Item(10, 17).Value = -1.25
$ws.Cells.Item(10, 18).Value = 1.85
$ws.Cells.Item(10, 19).Value = 1.95
$ws.Cells.Item(10, 20).Value = 2.75
$ws.Cells.Item(10, 21).Value = 1.8
$ws.Cells.Item(10, 22).Value = 2
$ws.Cells.Item(10, 23).Value = 0.3999999999999999
$ws.Cells.Item(10, 24).Value = -1
$ws.Cells.Item(10, 27).Value = 0.475
$ws.Cells.Item(10, 28).Value = -1
$ws.Cells.Item(10, 29).Value = 1
# Row 12
$ws.Cells.Item(12, 7).Value = "NK Domzale"
# Row 14
$ws.Cells.Item(14, 6).Value = "NK Maribor"
# Row 20
$ws.Cells.Item(20, 7).Value = "NK Maribor"
# Row 21
$ws.Cells.Item(21, 6).Value = "NK Domzale"
# Row 23
$ws.Cells.Item(23, 6).Value = "NK Maribor"
$ws.Cells.Item(23, 7).Value = "NK Domzale"
# Row 27
$ws.Cells.Item(27, 6).Value = "NK Domzale"
# Row 29
$ws.Cells.Item(29, 7).Value = "NK Maribor"
# Row 33
$ws.Cells.Item(33, 7).Value = "NK Maribor"
# Row 34
$ws.Cells.Item(34, 7).Value = "NK Domzale"
# Row 36
$ws.Cells.Item(36, 6).Value = "NK Domzale"
# Row 40
$ws.Cells.Item(40, 6).Value = "NK Maribor"
# Row 43
$ws.Cells.Item(43, 6).Value = "NK Domzale"
# Row 44
$ws.Cells.Item(44, 6).Value = "NK Maribor"
# Row 46
$ws.Cells.Item(46, 7).Value = "NK Domzale"
# Row 50
$ws.Cells.Item(50, 7).Value = "NK Maribor"
# Row 52
$ws.Cells.Item(52, 6).Value = "NK Maribor"
# Row 53
$ws.Cells.Item(53, 6).Value = "NK Domzale"
# Row 57
$ws.Cells.Item(57, 6).Value = "NK Domzale"
# Row 58
$ws.Cells.Item(58, 7).Value = "NK Maribor"
# Row 61
$ws.Cells.Item(61, 7).Value = "NK Domzale"
# Row 63
$ws.Cells.Item(63, 6).Value = "NK Maribor"
# Row 68
$ws.Cells.Item(68, 6).Value = "NK Domzale"
$ws.Cells.Item(68, 7).Value = "NK Maribor"
# Row 73
$ws.Cells.Item(73, 7).Value = "NK Domzale"
# Row 74
$ws.Cells.Item(74, 6).Value = "NK Maribor"
# Row 76
$ws.Cells.Item(76, 7).Value = "NK Domzale"
# Row 79
$ws.Cells.Item(79, 7).Value = "NK Maribor"
# Row 81
$ws.Cells.Item(81, 6).Value = "NK Domzale"
# Row 83
$ws.Cells.Item(83, 6).Value = "NK Maribor"
# Row 85
$ws.Cells.Item(85, 7).Value = "NK Domzale"
# Row 86
$ws.Cells.Item(86, 7).Value = "NK Maribor"
# Row 92
$ws.Cells.Item(92, 7).Value = "NK Domzale"
# Row 93
$ws.Cells.Item(93, 7).Value = "NK Maribor"
# Row 96
$ws.Cells.Item(96, 6).Value = "NK Domzale"
# Row 99
$ws.Cells.Item(99, 6).Value = "NK Maribor"
# Row 101
$ws.Cells.Item(101, 7).Value = "NK Domzale"
# Row 105
$ws.Cells.Item(105, 7).Value = "NK Domzale"
# Row 107
$ws.Cells.Item(107, 6).Value = "NK Maribor"
# Row 109
$ws.Cells.Item(109, 7).Value = "NK Maribor"
# Row 110
$ws.Cells.Item(110, 6).Value = "NK Domzale"
# Row 115
$ws.Cells.Item(115, 6).Value = "NK Maribor"
$ws.Cells.Item(115, 7).Value = "NK Domzale"
# Row 119
$ws.Cells.Item(119, 7).Value = "NK Maribor"
# Row 124
$ws.Cells.Item(124, 7).Value = "NK Domzale"
# Row 127
$ws.Cells.Item(127, 6).Value = "NK Domzale"
# Row 130
$ws.Cells.Item(130, 6).Value = "NK Maribor"
# Row 132
$ws.Cells.Item(132, 2).Value = 7977922
$ws.Cells.Item(132, 5).Value = 45380.5625
$ws.Cells.Item(132, 6).Value = "NK Maribor"
$ws.Cells.Item(132, 7).Value = "NK Radomlje"
$ws.Cells.Item(132, 11).Value = 1.285
$ws.Cells.Item(132, 12).Value = 5.5
$ws.Cells.Item(132, 13).Value = 6.5
$ws.Cells.Item(132, 14).Value = 1.4
$ws.Cells.Item(132, 15).Value = 5
$ws.Cells.Item(132, 16).Value = 5
$ws.Cells.Item(132, 17).Value = -1.25
$ws.Cells.Item(132, 18).Value = 1.975
$ws.Cells.Item(132, 19).Value = 1.825
$ws.Cells.Item(132, 20).Value = 2.75
$ws.Cells.Item(132, 21).Value = 1.875
$ws.Cells.Item(132, 22).Value = 1.925
# Row 133
$ws.Cells.Item(133, 2).Value = 7977924
$ws.Cells.Item(133, 5).Value = 45380.67708333334
$ws.Cells.Item(133, 7).Value = "NS Mura"
$ws.Cells.Item(133, 11).Value = 2
$ws.Cells.Item(133, 12).Value = 3.4
$ws.Cells.Item(133, 13).Value = 3.1
$ws.Cells.Item(133, 14).Value = 2.3
$ws.Cells.Item(133, 15).Value = 3.3
$ws.Cells.Item(133, 16).Value = 2.625
$ws.Cells.Item(133, 17).Value = 0
$ws.Cells.Item(133, 18).Value = 1.75
$ws.Cells.Item(133, 19).Value = 2.05
$ws.Cells.Item(133, 20).Value = 2.25
$ws.Cells.Item(133, 21).Value = 1.775
$ws.Cells.Item(133, 22).Value = 2.025
# Row 134
$ws.Cells.Item(134, 2).Value = 7977921
$ws.Cells.Item(134, 5).Value = 45381.45833333334
$ws.Cells.Item(134, 6).Value = "Olimpija Ljubljana"
$ws.Cells.Item(134, 7).Value = "FC Koper"
$ws.Cells.Item(134, 11).Value = 1.5
$ws.Cells.Item(134, 12).Value = 3.6
$ws.Cells.Item(134, 13).Value = 6
$ws.Cells.Item(134, 14).Value = 1.5
$ws.Cells.Item(134, 15).Value = 3.75
$ws.Cells.Item(134, 16).Value = 5.75
$ws.Cells.Item(134, 17).Value = -1
$ws.Cells.Item(134, 18).Value = 1.85
$ws.Cells.Item(134, 19).Value = 1.95
$ws.Cells.Item(134, 20).Value = 2.75
# Row 135
$ws.Cells.Item(135, 2).Value = 7977923
$ws.Cells.Item(135, 5).Value = 45381.5625
$ws.Cells.Item(135, 6).Value = "NK Celje"
$ws.Cells.Item(135, 7).Value = "NK Aluminij"
$ws.Cells.Item(135, 11).Value = 1.2
$ws.Cells.Item(135, 12).Value = 6.5
$ws.Cells.Item(135, 13).Value = 8
$ws.Cells.Item(135, 14).Value = 1.222
$ws.Cells.Item(135, 15).Value = 6.5
$ws.Cells.Item(135, 16).Value = 10
$ws.Cells.Item(135, 17).Value = -1.75
$ws.Cells.Item(135, 18).Value = 1.8
$ws.Cells.Item(135, 19).Value = 2
$ws.Cells.Item(135, 20).Value = 3
$ws.Cells.Item(135, 21).Value = 1.825
$ws.Cells.Item(135, 22).Value = 1.975

# Delete row 136 entirely (shifts dimension to AC135)
$ws.Rows.Item(136).Delete()
